$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the data rows 2, 4 and 5
# (row 3 is left untouched):
#   old row2 -> row4
#   old row4 -> row5
#   old row5 -> row2
#
# Capture the current ("before") values for the columns that change
# (D, J, K, L, M, N, O, P, Q) before overwriting anything.

$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

$row2 = @{}
$row4 = @{}
$row5 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("$col" + "2").Value2
    $row4[$col] = $ws.Range("$col" + "4").Value2
    $row5[$col] = $ws.Range("$col" + "5").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col" + "4").Value2 = $row2[$col]
    $ws.Range("$col" + "5").Value2 = $row4[$col]
    $ws.Range("$col" + "2").Value2 = $row5[$col]
}
